# Update "想去人数" (want-to-go count) figures in column F across all sheets.
# Generated from upstream data refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1558
$ws.Range("F3").Value = 3310
$ws.Range("F4").Value = 25
$ws.Range("F5").Value = 743
$ws.Range("F6").Value = 2338
$ws.Range("F7").Value = 501
$ws.Range("F8").Value = 420
$ws.Range("F9").Value = 251
$ws.Range("F11").Value = 364
$ws.Range("F12").Value = 1108
$ws.Range("F13").Value = 459
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 91
$ws.Range("F16").Value = 264
$ws.Range("F17").Value = 4809
$ws.Range("F19").Value = 1367
$ws.Range("F20").Value = 3555
$ws.Range("F21").Value = 338
$ws.Range("F22").Value = 140
$ws.Range("F24").Value = 3814
$ws.Range("F25").Value = 5190
$ws.Range("F27").Value = 983
$ws.Range("F28").Value = 571
$ws.Range("F29").Value = 3340
$ws.Range("F30").Value = 385
$ws.Range("F31").Value = 52
$ws.Range("F32").Value = 146
$ws.Range("F34").Value = 897
$ws.Range("F36").Value = 23
$ws.Range("F37").Value = 34
$ws.Range("F38").Value = 1432
$ws.Range("F39").Value = 142
$ws.Range("F40").Value = 1416
$ws.Range("F41").Value = 908
$ws.Range("F42").Value = 879
$ws.Range("F43").Value = 521
$ws.Range("F44").Value = 61
$ws.Range("F45").Value = 485
$ws.Range("F46").Value = 83
$ws.Range("F47").Value = 176
$ws.Range("F48").Value = 372
$ws.Range("F49").Value = 3748

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 1024
$ws.Range("F23").Value = 41

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2387

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2387
$ws.Range("F3").Value = 1558
$ws.Range("F4").Value = 3310
$ws.Range("F5").Value = 25
$ws.Range("F6").Value = 743
$ws.Range("F7").Value = 2338
$ws.Range("F8").Value = 501
$ws.Range("F9").Value = 420
$ws.Range("F10").Value = 251
$ws.Range("F11").Value = 1024
$ws.Range("F13").Value = 364
$ws.Range("F14").Value = 1108
$ws.Range("F15").Value = 459
$ws.Range("F16").Value = 152
$ws.Range("F17").Value = 91
$ws.Range("F18").Value = 264
$ws.Range("F19").Value = 4809
$ws.Range("F20").Value = 1367
$ws.Range("F21").Value = 3814
$ws.Range("F22").Value = 5190
$ws.Range("F24").Value = 983
$ws.Range("F25").Value = 571
$ws.Range("F26").Value = 3340
$ws.Range("F27").Value = 385
$ws.Range("F28").Value = 52
$ws.Range("F29").Value = 146
$ws.Range("F31").Value = 897
$ws.Range("F33").Value = 23
$ws.Range("F34").Value = 34
$ws.Range("F35").Value = 1432
$ws.Range("F36").Value = 1416
$ws.Range("F37").Value = 908
$ws.Range("F39").Value = 521
$ws.Range("F41").Value = 61
$ws.Range("F42").Value = 41
$ws.Range("F43").Value = 490
$ws.Range("F45").Value = 83
$ws.Range("F46").Value = 176
$ws.Range("F47").Value = 372
$ws.Range("F49").Value = 3748
